$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '71.065.29'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.806.78'
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '701.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '3.805.71'
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("D15").Value = '4.450.95'
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("D16").Value = '3.809.89'
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("D17").Value = '71.104.81'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '512.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.05%  '
$ws.Range("E23").Value = '  -0.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("D26").Value = '3.958.22'
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.41%  '
$ws.Range("E31").Value = '  -4.96%  '
$ws.Range("E32").Value = '  -1.36%  '
$ws.Range("E33").Value = '  -1.60%  '
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = '3.770.14'
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '172.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.95%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '429.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.37%  '
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.61'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.11%  '
